# Applies the "output generated at 456a3b4" update to the 合肥-漫展信息 workbook.
#
# For both the "展览" sheet and the "全部类型" sheet, a new event row
# (合肥·原神X星铁Only, 2024-06-29) is inserted right before the existing
# "合肥·第1.5届星芒动漫嘉年华" (2024-06-30) row, pushing every row below it
# down by one. A handful of unrelated "想去人数" (F column) counters also
# ticked up slightly between scrapes.

$wb = $excel.ActiveWorkbook

function Update-MaoZhanSheet {
    param(
        [string]$SheetName
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Locate the anchor row ("合肥·第1.5届星芒动漫嘉年华", 2024-06-30) by its
    # name in column C so this works regardless of the sheet's exact layout.
    $anchorRow = $ws.Columns.Item(3).Find("合肥·第1.5届星芒动漫嘉年华").Row

    # A few "想去人数" counters (column F) ticked up slightly since the last
    # scrape. Find + bump each one before the insert shifts anything. (The
    # row number is pulled into a plain variable before use — chaining
    # straight off Find(...).EntireRow does not reliably target that row.)
    $r = $ws.Columns.Item(3).Find("合肥·运动番only·群青日和").Row
    $ws.Cells.Item($r, 6).Value = 651
    $r = $ws.Columns.Item(3).Find("合肥·FT动漫嘉年华（免费）").Row
    $ws.Cells.Item($r, 6).Value = 226
    $r = $ws.Columns.Item(3).Find("合肥·第六届环形宇宙动漫游戏嘉年华-一周年超强巨制~").Row
    $ws.Cells.Item($r, 6).Value = 9888
    $r = $ws.Columns.Item(3).Find("合肥·城市动漫节").Row
    $ws.Cells.Item($r, 6).Value = 3951
    $r = $ws.Columns.Item(3).Find("合肥·第七届环形宇宙动漫游戏嘉年华").Row
    $ws.Cells.Item($r, 6).Value = 1465

    # Insert a blank row right above the anchor row; everything from the
    # anchor row down shifts one row further down.
    $ws.Rows.Item($anchorRow).Insert()

    # Fill in the newly inserted row with the new event's data.
    $newRow = $anchorRow

    $ws.Cells.Item($newRow, 1).Value = $anchorRow - 1

    # Column B holds plain date-like text ("2024-06-29"); force text
    # formatting first so Excel doesn't silently convert it to a date serial.
    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = "2024-06-29"

    $ws.Cells.Item($newRow, 3).Value = "合肥·原神X星铁Only"
    $ws.Cells.Item($newRow, 4).Value = "金寨路与天堂窄路交叉口 梵木艺术中心"
    $ws.Cells.Item($newRow, 5).Value = "2024.06.29 10:00-06.29 17:00"
    $ws.Cells.Item($newRow, 6).Value = 0
    $ws.Cells.Item($newRow, 7).Value = 60
    $ws.Cells.Item($newRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86406"
    $ws.Cells.Item($newRow, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/r3c5IueN1716820859877.jpeg"

    # Renumber the serial-number column (A) for every row that shifted down.
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = $newRow + 1; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

Update-MaoZhanSheet "展览"
Update-MaoZhanSheet "全部类型"
